$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.501.33"
$ws.Range("E2").Value = "  +12.12%  "

$ws.Range("D3").Value = "1.832.52"
$ws.Range("E3").Value = "  +8.14%  "

$ws.Range("E4").Value = "  +0.19%  "

$cell = $ws.Range("D5")
$cell.Value = "'231.66"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +4.46%  "

$cell = $ws.Range("D6")
$cell.Value = "'0.547"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +4.48%  "

$ws.Range("E7").Value = "  +0.28%  "

$cell = $ws.Range("D8")
$cell.Value = "'31.66"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +1.92%  "

$ws.Range("E9").Value = "  +1.18%  "

$ws.Range("E10").Value = "  +6.24%  "

$ws.Range("E11").Value = "  +8.17%  "

$ws.Range("E12").Value = "  +3.41%  "

$ws.Range("D13").Value = "2.092.38"
$ws.Range("E13").Value = "  +8.05%  "

$ws.Range("D14").Value = "1.838.75"
$ws.Range("E14").Value = "  +8.46%  "

$ws.Range("E15").Value = "  +4.79%  "

$ws.Range("D16").Value = "34.454.76"
$ws.Range("E16").Value = "  +11.93%  "

$cell = $ws.Range("D17")
$cell.Value = "'10.33"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -3.88%  "

$cell = $ws.Range("D18")
$cell.Value = "'4.38"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +8.79%  "

$ws.Range("E19").Value = "  +5.72%  "

$cell = $ws.Range("D20")
$cell.Value = "'260.78"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +4.84%  "

$ws.Range("D21").Value = "0.0₃0755"
$ws.Range("E21").Value = "  +4.74%  "

$cell = $ws.Range("D22")
$cell.Value = "'0.998"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.13%  "

$cell = $ws.Range("D23")
$cell.Value = "'10.57"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +3.30%  "

$cell = $ws.Range("D24")
$cell.Value = "'4.41"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +2.30%  "

$ws.Range("E25").Value = "  +0.77%  "

$cell = $ws.Range("D26")
$cell.Value = "'161.18"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +2.33%  "

$cell = $ws.Range("D27")
$cell.Value = "'16.88"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +5.60%  "

$cell = $ws.Range("D28")
$cell.Value = "'7.22"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +6.68%  "

$ws.Range("E29").Value = "  +4.62%  "

$ws.Range("E30").Value = "  +0.17%  "

$cell = $ws.Range("D31")
$cell.Value = "'3.86"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +10.39%  "

$cell = $ws.Range("D32")
$cell.Value = "'0.0519"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +3.53%  "

$ws.Range("E33").Value = "  +7.23%  "

$ws.Range("E34").Value = "  +8.56%  "

$ws.Range("D35").Value = "1.583.23"
$ws.Range("E35").Value = "  +4.07%  "

$ws.Range("E36").Value = "  +6.00%  "

$cell = $ws.Range("D37")
$cell.Value = "'1.08"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +4.18%  "

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell = $ws.Range("D38")
$cell.Value = "'0.636"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +8.49%  "

$ws.Range("B39").Value = "Aave"
$ws.Range("C39").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$cell = $ws.Range("D39")
$cell.Value = "'85.65"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +7.16%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cell = $ws.Range("D40")
$cell.Value = "'0.0190"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +5.54%  "

$ws.Range("E41").Value = "  +5.41%  "

$ws.Range("E42").Value = "  +1.64%  "

$cell = $ws.Range("D43")
$cell.Value = "'0.923"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +7.64%  "

$cell = $ws.Range("D44")
$cell.Value = "'2.15"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +6.19%  "

$cell = $ws.Range("D45")
$cell.Value = "'0.0522"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +3.74%  "

$cell = $ws.Range("D46")
$cell.Value = "'1.07"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +4.56%  "

$ws.Range("D47").Value = "1.983.61"
$ws.Range("E47").Value = "  +8.21%  "

$cell = $ws.Range("D48")
$cell.Value = "'5.76"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +6.02%  "

$cell = $ws.Range("D49")
$cell.Value = "'53.28"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +1.34%  "

$ws.Range("E50").Value = "  +0.24%  "

$ws.Range("E51").Value = "  +9.88%  "
